$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Biden for President count updated 11 -> 18 (row 45)
$ws.Range("B45").Value = 18

# New entry "North Carolina Democratic Party" (n=1) inserted alphabetically
# right before "Northeastern University" (currently row 269), shifting all
# subsequent rows down by one.
$ws.Rows("269:269").Insert()
$ws.Range("A269").Value = "North Carolina Democratic Party"
$ws.Range("B269").Value = 1
